$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "69.564.80"
Set-TextValue "E2" "  +2.05%  "
Set-TextValue "D3" "3.381.65"
Set-TextValue "E3" "  +1.19%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "581.34"
Set-TextValue "E5" "  -0.44%  "
Set-TextValue "D6" "179.00"
Set-TextValue "E6" "  +0.97%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.05%  "
Set-TextValue "E8" "  +0.69%  "
Set-TextValue "D9" "0.198"
Set-TextValue "E9" "  +8.43%  "
Set-TextValue "E10" "  +1.09%  "
Set-TextValue "D11" "48.51"
Set-TextValue "E11" "  +1.00%  "
Set-TextValue "E12" "  +4.16%  "
Set-TextValue "D13" "688.61"
Set-TextValue "E13" "  -0.72%  "
Set-TextValue "D14" "8.63"
Set-TextValue "E14" "  +2.54%  "
Set-TextValue "D15" "3.923.64"
Set-TextValue "E15" "  +1.04%  "
Set-TextValue "D16" "69.563.54"
Set-TextValue "E16" "  +1.91%  "
Set-TextValue "E17" "  +0.87%  "
Set-TextValue "D18" "3.375.65"
Set-TextValue "E18" "  +0.97%  "
Set-TextValue "D19" "17.73"
Set-TextValue "E19" "  +1.82%  "
Set-TextValue "E20" "  +0.67%  "
Set-TextValue "E21" "  +1.79%  "
Set-TextValue "D22" "17.28"
Set-TextValue "E22" "  +1.66%  "
Set-TextValue "D23" "5.35"
Set-TextValue "E23" "  -1.70%  "
Set-TextValue "D24" "101.66"
Set-TextValue "E24" "  +1.45%  "
Set-TextValue "D27" "9.69"
Set-TextValue "E27" "  +1.71%  "
Set-TextValue "D28" "33.62"
Set-TextValue "E28" "  +2.00%  "
Set-TextValue "D29" "8.73"
Set-TextValue "E29" "  +2.71%  "
Set-TextValue "D30" "6.89"
Set-TextValue "E30" "  -0.85%  "
Set-TextValue "D31" "3.87"
Set-TextValue "E31" "  +18.55%  "
Set-TextValue "E32" "  +0.24%  "
Set-TextValue "D33" "554.56"
Set-TextValue "E33" "  -1.98%  "
Set-TextValue "E34" "  +0.34%  "
Set-TextValue "D35" "58.03"
Set-TextValue "E35" "  +1.00%  "
Set-TextValue "E36" "  +0.14%  "
Set-TextValue "D37" "3.603.19"
Set-TextValue "E37" "  -2.35%  "
Set-TextValue "E38" "  +2.72%  "
Set-TextValue "D39" "35.31"
Set-TextValue "E39" "  +1.67%  "
Set-TextValue "D40" "0.0₃0727"
Set-TextValue "E40" "  +8.44%  "
Set-TextValue "D41" "3.31"
Set-TextValue "E41" "  +3.98%  "
Set-TextValue "E42" "  +4.72%  "
Set-TextValue "D43" "3.38"
Set-TextValue "E43" "  +2.70%  "
Set-TextValue "E44" "  +2.90%  "
Set-TextValue "E45" "  +0.24%  "
Set-TextValue "E46" "  +0.29%  "
Set-TextValue "E47" "  +0.43%  "
Set-TextValue "B48" "FirstDigitalUSD"
Set-TextValue "C48" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  -0.26%  "
Set-TextValue "B49" "Mantle"
Set-TextValue "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "1.38"
Set-TextValue "E49" "  +3.55%  "
Set-TextValue "D50" "129.12"
Set-TextValue "E50" "  -1.39%  "
Set-TextValue "D51" "2.56"
Set-TextValue "E51" "  -0.48%  "

Write-Host "Applied all changes"
